{"js": "// Remove hard line-break artifacts that were baked into a few <w:t> runs:\n// - trailing \"\\n\" at the end of a run's text is dropped\n// - embedded \"\\n\" line-wrap characters inside a run's text become a single space\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map of exact original paragraph text -> desired replacement text.\nconst replacements = new Map([\n  [\n    \"This is an annotatable resource in the casebook.\\n\",\n    \"This is an annotatable resource in the casebook.\",\n  ],\n  [\n    \"highlighted: content to highlight; elided: content to elide;\\nreplaced: content to replace; linked: content to link; noted:\\ncontent to note; highlighted2: second highlight content;\",\n    \"highlighted: content to highlight; elided: content to elide; replaced: content to replace; linked: content to link; noted: content to note; highlighted2: second highlight content;\",\n  ],\n  [\n    \"This is the second chapter of the casebook.\\n\",\n    \"This is the second chapter of the casebook.\",\n  ],\n]);\n\nfor (const paragraph of paragraphs.items) {\n  const newText = replacements.get(paragraph.text);\n  if (newText !== undefined) {\n    paragraph.getRange().insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove hard line-break artifacts that were baked into a few runs' text:\n# - a trailing LF at the end of a run's text is dropped\n# - an embedded LF line-wrap character inside a run's text becomes a single space\n$d = $word.ActiveDocument\n$nl = [char]10\n\n$replacements = @(\n    @{\n        Find    = \"This is an annotatable resource in the casebook.${nl}\"\n        Replace = \"This is an annotatable resource in the casebook.\"\n    },\n    @{\n        Find    = \"highlighted: content to highlight; elided: content to elide;${nl}replaced: content to replace; linked: content to link; noted:${nl}content to note; highlighted2: second highlight content;\"\n        Replace = \"highlighted: content to highlight; elided: content to elide; replaced: content to replace; linked: content to link; noted: content to note; highlighted2: second highlight content;\"\n    },\n    @{\n        Find    = \"This is the second chapter of the casebook.${nl}\"\n        Replace = \"This is the second chapter of the casebook.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Find\n    $find.Replacement.Text = $r.Replace\n    $find.Execute($r.Find, $false, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2) | Out-Null\n}\n"}
